# adding averages and more checks
#
# Applies:
#  1) Training Dashboard: PERIOD TO EXPIRE (col H) decreases by 8 for rows 3-15,
#     and LAST UPDATE (col I) moves from 08-Sep-2025 to 16-Sep-2025 for rows 3-15.
#  2) Exam Dashboard: COMMENTS column (E) text changes from "OK" to "date is valid"
#     for rows 3-9, and column E is widened.
#  3) Header row / title styling: bold white font (was bold black for headers,
#     bold size-14 default-color for the titles).

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1) Training Dashboard - PERIOD TO EXPIRE (H) & LAST UPDATE (I)
# ---------------------------------------------------------------------------
$trainingWs = $wb.Worksheets.Item("Training Dashboard")

$periodToExpire = @{3=523;4=521;5=523;6=521;7=523;8=525;9=532;10=525;11=539;12=609;13=323;14=268;15=323}

foreach ($row in $periodToExpire.Keys) {
    $trainingWs.Cells.Item($row, 8).Value = $periodToExpire[$row]
}

# Force the LAST UPDATE column to stay text (it was stored as literal text in
# the workbook, not a real date serial) before writing the new date string.
$lastUpdateRange = $trainingWs.Range("I3:I15")
$lastUpdateRange.NumberFormat = "@"
for ($row = 3; $row -le 15; $row++) {
    $trainingWs.Cells.Item($row, 9).Value = "16-Sep-2025"
}

# ---------------------------------------------------------------------------
# 2) Exam Dashboard - COMMENTS (E) text + column width
# ---------------------------------------------------------------------------
$examWs = $wb.Worksheets.Item("Exam Dashboard")

for ($row = 3; $row -le 9; $row++) {
    $examWs.Cells.Item($row, 5).Value = "date is valid"
}

# Widen column E from 10 to 15 (OOXML character-width units); ColumnWidth is
# offset from the saved <col width> by 5/6 of a character.
$examWs.Columns.Item(5).ColumnWidth = 15 - 5/6

# ---------------------------------------------------------------------------
# 3) Header / title styling - bold white font
# ---------------------------------------------------------------------------
foreach ($sheetName in @("Training Dashboard", "Exam Dashboard")) {
    $ws = $wb.Worksheets.Item($sheetName)

    # Title cell (row 1): drop the old 14pt size, make it bold white (same
    # font the header row now uses).
    $ws.Range("A1").Font.Size = 11
    $ws.Range("A1").Font.Color = 16777215

    # Header row (row 2): bold white text on its existing dark-blue fill.
    $ws.Rows.Item(2).Font.Color = 16777215
}
